$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B1").Value = "California"
$ws.Range("C1").Value = (Get-Date -Year 2022 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0)
